# Update countries & provincias Spain
# Applies the 28-Mar-2020 15:59 -> 16:29 data refresh to the "Pais" sheet:
#  - timestamp banner text (A1)
#  - updated case counts for several countries (B:H)
#  - re-ordering of a few countries that changed rank when sorted by total cases
#    (the country label in column A moves to a different row, carrying its
#    own new numbers with it)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "last refreshed" banner -------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 16:29"

# --- Helper: write a full row (Country name + 7 numeric columns) ---------
function Set-CountryRow {
    param($Row, $Country, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes)
    $ws.Cells.Item($Row, 1).Value = $Country
    $ws.Cells.Item($Row, 2).Value = $Total
    $ws.Cells.Item($Row, 3).Value = $Nuevos
    $ws.Cells.Item($Row, 4).Value = $Activos
    $ws.Cells.Item($Row, 5).Value = $Recuperados
    $ws.Cells.Item($Row, 6).Value = $Criticos
    $ws.Cells.Item($Row, 7).Value = $MuertesHoy
    $ws.Cells.Item($Row, 8).Value = $Muertes
}

# Row 4: Estados Unidos - numbers only
Set-CountryRow 4 "Estados Unidos" 105161 1035 2538 100901 2494 26 1722

# Row 7: Espana - numbers only
Set-CountryRow 7 "España" 72248 6529 12285 54151 4165 674 5812

# Rows 29-30: Chile/Ecuador swap rank and refresh numbers
Set-CountryRow 29 "Chile" 1909 299 43 1861 7 0 5
Set-CountryRow 30 "Ecuador" 1627 0 3 1583 58 0 41

# Row 40: Finlandia - numbers only
$ws.Cells.Item(40, 5).Value = 1148
$ws.Cells.Item(40, 7).Value = 2
$ws.Cells.Item(40, 8).Value = 9

# Rows 54-56: Croacia/Estonia/Peru shift rank and refresh numbers
Set-CountryRow 54 "Croacia" 657 71 45 607 14 2 5
Set-CountryRow 55 "Estonia" 645 70 20 624 10 0 1
Set-CountryRow 56 "Peru" 635 0 16 608 21 0 11

# Rows 73-77: Letonia/Eslovaquia/Taiwan/Uruguay/Principado de Andorra shift rank
Set-CountryRow 73 "Principado de Andorra" 308 41 1 304 10 0 3
Set-CountryRow 74 "Letonia" 305 25 1 304 3 0 0
Set-CountryRow 75 "Eslovaquia" 292 23 2 290 1 0 0
Set-CountryRow 76 "Taiwan" 283 16 30 251 0 0 2
Set-CountryRow 77 "Uruguay" 274 36 0 274 8 0 0

# Row 117: Trinidad y Tobago - numbers only
$ws.Cells.Item(117, 2).Value = 74
$ws.Cells.Item(117, 3).Value = 8
$ws.Cells.Item(117, 5).Value = 71
